# ---------------------------------------------------------------------------
# "General edits on shell presentations and added my name"
#
# 1) Bump the cached datetimeFigureOut field text on the two title-slide
#    layouts (Title Slide - Portrait image / Title Slide - Landscape image)
#    from 03/10/2018 -> 12/10/2018.
# 2) Slide 1 ("The Unix Shell"): merge the two title runs ("The " + "Unix
#    Shell") into a single run, and drop the stray trailing endParaRPr on
#    the subtitle ("Text Editors and Terminals").
# 3) Slide 2 ("What is an Editor"): fix the "An text editor" typo, and
#    re-word "...save them in a file." -> "...save them in a file" (drop
#    the trailing full stop), which ends up as three runs.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Date placeholder on the two affected slide layouts -----------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "03/10/2018") {
                $tr.Text = "12/10/2018"
            }
        }
    }
}

# --- 2) Slide 1 : title + subtitle -----------------------------------------
$slide1 = $p.Slides.Item(1)

$title = $slide1.Shapes.Item(1)
# Force a genuine rewrite (route through an unrelated placeholder first) so
# the two existing runs get collapsed into a single run instead of the
# no-op / partial-diff path that a same-looking string would take.
$title.TextFrame.TextRange.Text = "ZZZZZZZZZZ"
$title.TextFrame.TextRange.Text = "The Unix Shell"

$subtitle = $slide1.Shapes.Item(2)
# Drop the dangling endParaRPr: clearing the frame removes it, then we can
# retype the text without it coming back.
$subtitle.TextFrame.TextRange.Delete()
$subtitle.TextFrame.TextRange.Text = "Text Editors and Terminals"

# --- 3) Slide 2 : body copy -------------------------------------------------
$slide2 = $p.Slides.Item(2)
$body = $slide2.Shapes.Item(2)
$bodyRange = $body.TextFrame.TextRange

# Paragraph 1: "An text editor..." -> "A text editor..."
$para1 = $bodyRange.Paragraphs(1, 1)
$para1.Text = "ZZZZZZZZZZ"
$body.TextFrame.TextRange.Paragraphs(1, 1).Text = "A text editor is a program used to edit text!"

# Paragraph 2: drop the trailing full stop and re-type "a file" as its own
# runs (mirrors the authored edit, which split the sentence into three
# runs: the unchanged lead-in, "a ", and "file").
$para2 = $body.TextFrame.TextRange.Paragraphs(2, 1)
$leadIn = "There are lots of them, but they all manipulate a stream of characters so you can save them in "
$oldTail = $para2.Characters($para2.Length - 7, 7)
$oldTail.Text = ""
$para2.InsertAfter("a ")
$para2.InsertAfter("file")
